# Applies the M2Doc "field -> plain text" rewrite described in the commit:
#   "Updated parser to use TokenIteratorFieldRewriterSplit."
#
# 1) Remove the stray 4-space run that sits right before "demonstration"
#    in the first paragraph (the one between "...NotExistingRepo" and
#    "demonstration").
# 2) Replace the `{ m:self.name }` Word field (fldChar begin/instrText/end)
#    in the second paragraph with plain literal text runs "{", "m", ":",
#    "self" (keeping its orange theme colour) and ".name}".

$d = $word.ActiveDocument

# --- Step 1: delete the duplicate "    " run before "demonstration" -------
$matchRange = $d.Content
$found = $matchRange.Find.Execute("NotExistingRepo    demonstration")
if (-not $found) {
    throw "Could not locate 'NotExistingRepo    demonstration' text"
}
$spacesStart = $matchRange.Start + ("NotExistingRepo").Length
$spacesEnd = $spacesStart + 4
$spacesRange = $d.Range($spacesStart, $spacesEnd)
if ($spacesRange.Text -ne "    ") {
    throw "Unexpected text when deleting spaces run: [" + $spacesRange.Text + "]"
}
$spacesRange.Delete()

# --- Step 2: turn the `{ m:self.name }` field into literal text runs ------
if ($d.Fields.Count -ne 1) {
    throw "Expected exactly one field, found " + $d.Fields.Count
}
$fieldParagraph = $d.Paragraphs.Item(2)
$paragraphStart = $fieldParagraph.Range.Start
$field = $d.Fields.Item(1)
$field.Delete()
$target = $d.Range($paragraphStart, $paragraphStart)

$newParaXml = '<w:body>' +
    '<w:p>' +
        '<w:r><w:t>{</w:t></w:r>' +
        '<w:r><w:t>m</w:t></w:r>' +
        '<w:r><w:t>:</w:t></w:r>' +
        '<w:r><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t>self</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">.name}</w:t></w:r>' +
    '</w:p>' +
'</w:body>'

$packageXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
            '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                $newParaXml +
            '</w:document>' +
        '</pkg:xmlData>' +
    '</pkg:part>' +
'</pkg:package>'

$target.InsertXML($packageXml)
